# Updated cryptos list with latest Price / Volume(1h) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dCell = $ws.Cells.Item(2, 4)
$dCell.Value = "'30.516.27"
$dCell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.09%  "

$dCell = $ws.Cells.Item(3, 4)
$dCell.Value = "'1.887.57"
$dCell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +0.79%  "

$ws.Cells.Item(4, 5).Value = "  -0.11%  "

$dCell = $ws.Cells.Item(5, 4)
$dCell.Value = "'244.05"
$dCell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.35%  "

$ws.Cells.Item(6, 5).Value = "  -0.06%  "

$dCell = $ws.Cells.Item(7, 4)
$dCell.Value = "'0.4720"
$dCell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.17%  "

$dCell = $ws.Cells.Item(8, 4)
$dCell.Value = "'0.2894"
$dCell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.30%  "

$dCell = $ws.Cells.Item(9, 4)
$dCell.Value = "'0.06485"
$dCell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.29%  "

$ws.Cells.Item(10, 5).Value = "  +1.26%  "

$dCell = $ws.Cells.Item(11, 4)
$dCell.Value = "'0.07760"
$dCell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.92%  "

$dCell = $ws.Cells.Item(12, 4)
$dCell.Value = "'1.888.11"
$dCell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.85%  "

$dCell = $ws.Cells.Item(13, 4)
$dCell.Value = "'95.88"
$dCell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.07%  "

$dCell = $ws.Cells.Item(14, 4)
$dCell.Value = "'0.7264"
$dCell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.80%  "

$ws.Cells.Item(15, 5).Value = "  +0.57%  "

$dCell = $ws.Cells.Item(16, 4)
$dCell.Value = "'281.96"
$dCell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +3.17%  "

$dCell = $ws.Cells.Item(17, 4)
$dCell.Value = "'30.513.81"
$dCell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.07%  "

$dCell = $ws.Cells.Item(18, 4)
$dCell.Value = "'13.06"
$dCell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -1.35%  "

$ws.Cells.Item(19, 5).Value = "  -0.05%  "

$dCell = $ws.Cells.Item(20, 4)
$dCell.Value = "'0.000007475"
$dCell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.14%  "

$dCell = $ws.Cells.Item(21, 4)
$dCell.Value = "'2.136.98"
$dCell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.94%  "

$ws.Cells.Item(22, 5).Value = "  -0.16%  "

$ws.Cells.Item(23, 5).Value = "  +0.59%  "

$dCell = $ws.Cells.Item(24, 4)
$dCell.Value = "'6.334"
$dCell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +2.72%  "

$ws.Cells.Item(25, 5).Value = "  +0.00%  "

$dCell = $ws.Cells.Item(26, 4)
$dCell.Value = "'9.079"
$dCell.Style = "Normal"

$dCell = $ws.Cells.Item(27, 4)
$dCell.Value = "'18.87"
$dCell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +1.04%  "

$dCell = $ws.Cells.Item(28, 4)
$dCell.Value = "'1.893"
$dCell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -0.41%  "

$ws.Cells.Item(29, 5).Value = "  -0.67%  "

$dCell = $ws.Cells.Item(30, 4)
$dCell.Value = "'0.09673"
$dCell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -2.89%  "

$dCell = $ws.Cells.Item(31, 4)
$dCell.Value = "'1.472"
$dCell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -2.47%  "

$dCell = $ws.Cells.Item(32, 4)
$dCell.Value = "'4.278"
$dCell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.14%  "

$dCell = $ws.Cells.Item(33, 4)
$dCell.Value = "'4.148"
$dCell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +1.70%  "

$dCell = $ws.Cells.Item(34, 4)
$dCell.Value = "'0.04863"
$dCell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +2.14%  "

$dCell = $ws.Cells.Item(35, 4)
$dCell.Value = "'1.125"
$dCell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.84%  "

$dCell = $ws.Cells.Item(36, 4)
$dCell.Value = "'0.6931"
$dCell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.22%  "

$dCell = $ws.Cells.Item(37, 4)
$dCell.Value = "'2.714"
$dCell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.14%  "

$dCell = $ws.Cells.Item(38, 4)
$dCell.Value = "'0.01885"
$dCell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +1.91%  "

$dCell = $ws.Cells.Item(39, 4)
$dCell.Value = "'2.823"
$dCell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +2.66%  "

$dCell = $ws.Cells.Item(40, 4)
$dCell.Value = "'74.84"
$dCell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +2.54%  "

$dCell = $ws.Cells.Item(41, 4)
$dCell.Value = "'6.204"
$dCell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.41%  "

$dCell = $ws.Cells.Item(42, 4)
$dCell.Value = "'1.967"
$dCell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.25%  "

$dCell = $ws.Cells.Item(43, 4)
$dCell.Value = "'0.4270"
$dCell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +2.98%  "

$dCell = $ws.Cells.Item(44, 4)
$dCell.Value = "'0.9998"
$dCell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.12%  "

$dCell = $ws.Cells.Item(45, 4)
$dCell.Value = "'0.8258"
$dCell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.80%  "

$dCell = $ws.Cells.Item(46, 4)
$dCell.Value = "'101.27"
$dCell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.11%  "

$dCell = $ws.Cells.Item(47, 4)
$dCell.Value = "'9.636"
$dCell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +2.83%  "

$dCell = $ws.Cells.Item(48, 4)
$dCell.Value = "'6.956"
$dCell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.14%  "

$dCell = $ws.Cells.Item(49, 4)
$dCell.Value = "'35.17"
$dCell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.27%  "

$ws.Cells.Item(50, 5).Value = "  -0.45%  "

$dCell = $ws.Cells.Item(51, 4)
$dCell.Value = "'0.05747"
$dCell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +1.60%  "
